$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 107: Another Man's Ink / Enchanted Truegold Ink
$ws.Range("H107").Value = 873.4
$ws.Range("I107").Value = 941.75
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 941.75
$ws.Range("L107").Value = 600
$ws.Range("M107").Value = 978.25
$ws.Range("N107").Value = -4440
# Row 125: Body over Mind / Grade 5 Dexterity Alkahest
$ws.Range("H125").Value = 251524.75
$ws.Range("I125").Value = 599
$ws.Range("K125").Value = 5391
$ws.Range("M125").Value = -2931
# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 2762.6
$ws.Range("I132").Value = 2309.1562
$ws.Range("J132").Value = 7599.3335
$ws.Range("K132").Value = 6927.4686
$ws.Range("L132").Value = 22798.0005
$ws.Range("M132").Value = -4397.4686
$ws.Range("N132").Value = -27858.0005

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 16973.342
$ws.Range("I32").Value = 7500.1
$ws.Range("K32").Value = 7500.1
$ws.Range("M32").Value = -7213.1
# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 2610.5908
$ws.Range("J45").Value = 3754
$ws.Range("L45").Value = 3754
$ws.Range("N45").Value = -4508
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 2998.3845
$ws.Range("I61").Value = 1490.1875
$ws.Range("J61").Value = 5411.5
$ws.Range("K61").Value = 1490.1875
$ws.Range("L61").Value = 5411.5
$ws.Range("M61").Value = -1278.1875
$ws.Range("N61").Value = -5835.5
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 2840.926
$ws.Range("I132").Value = 2290.6365
$ws.Range("J132").Value = 5262.2
$ws.Range("K132").Value = 6871.9095
$ws.Range("L132").Value = 15786.6
$ws.Range("M132").Value = -4341.9095
$ws.Range("N132").Value = -20846.6
# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 2998.3845
$ws.Range("I136").Value = 1490.1875
$ws.Range("J136").Value = 5411.5
$ws.Range("K136").Value = 4470.5625
$ws.Range("L136").Value = 16234.5
$ws.Range("M136").Value = -1920.5625
$ws.Range("N136").Value = -21334.5
# Row 139: Backing up My Words / Titanium Gold Thornplate of Fending
$ws.Range("H139").Value = 80525.5
$ws.Range("J139").Value = 115496
$ws.Range("L139").Value = 115496
$ws.Range("N139").Value = -125776

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94: High Steal / High Steel Nugget
$ws.Range("H94").Value = 1233.7778
$ws.Range("I94").Value = 1259.3334
$ws.Range("K94").Value = 1259.3334
$ws.Range("M94").Value = -808.3334
# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 4486.647
$ws.Range("I105").Value = 3905.6924
$ws.Range("J105").Value = 6374.75
$ws.Range("K105").Value = 3905.6924
$ws.Range("L105").Value = 6374.75
$ws.Range("M105").Value = -2158.6924
$ws.Range("N105").Value = -9868.75
# Row 106: Fire for Hire / Molybdenum Rimfire
$ws.Range("H106").Value = 29000
$ws.Range("J106").Value = 29000
$ws.Range("L106").Value = 29000
$ws.Range("N106").Value = -31524
# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 6390.409
$ws.Range("I134").Value = 4400.353
$ws.Range("J134").Value = 13156.6
$ws.Range("K134").Value = 13201.059
$ws.Range("L134").Value = 39469.8
$ws.Range("M134").Value = -10666.059
$ws.Range("N134").Value = -44539.8
# Row 138: Bladewinner / Titanium Gold Greatsword
$ws.Range("H138").Value = 111999.4
$ws.Range("I138").Value = 49999
$ws.Range("K138").Value = 49999
$ws.Range("M138").Value = -44859

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall / Elm Lumber
$ws.Range("H22").Value = 1380.2
$ws.Range("J22").Value = 1145.8
$ws.Range("L22").Value = 1145.8
$ws.Range("N22").Value = -1845.8
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 9089.343000000001
$ws.Range("I31").Value = 3501.2
$ws.Range("J31").Value = 11324.6
$ws.Range("K31").Value = 3501.2
$ws.Range("L31").Value = 11324.6
$ws.Range("M31").Value = -3206.2
$ws.Range("N31").Value = -11914.6
# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 9089.343000000001
$ws.Range("I34").Value = 3501.2
$ws.Range("J34").Value = 11324.6
$ws.Range("K34").Value = 3501.2
$ws.Range("L34").Value = 11324.6
$ws.Range("M34").Value = -3299.2
$ws.Range("N34").Value = -11728.6
# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 4670.143
$ws.Range("I134").Value = 3670.12
$ws.Range("J134").Value = 13003.667
$ws.Range("K134").Value = 11010.36
$ws.Range("L134").Value = 39011.001
$ws.Range("M134").Value = -8475.360000000001
$ws.Range("N134").Value = -44081.001
# Row 135: The Wing's Wings / Ceiba Wings
$ws.Range("H135").Value = 80771
$ws.Range("J135").Value = 80771
$ws.Range("L135").Value = 80771
$ws.Range("N135").Value = -90911

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 47: Winter of Our Discontent / Mugwort Carp
$ws.Range("H47").Value = 474.33334
$ws.Range("I47").Value = 283.625
$ws.Range("J47").Value = 2000
$ws.Range("K47").Value = 850.875
$ws.Range("L47").Value = 6000
$ws.Range("M47").Value = -419.875
$ws.Range("N47").Value = -6862
# Row 56: Culture Club / Crowned Pie
$ws.Range("H56").Value = 3989
$ws.Range("I56").Value = 3989
$ws.Range("K56").Value = 3989
$ws.Range("M56").Value = -3459
# Row 106: Herky Jerky / Jerked Jhammel
$ws.Range("H106").Value = 10372.134
$ws.Range("I106").Value = 2666.6667
$ws.Range("J106").Value = 12298.5
$ws.Range("K106").Value = 8000.000100000001
$ws.Range("L106").Value = 36895.5
$ws.Range("M106").Value = -7054.000100000001
$ws.Range("N106").Value = -38787.5
# Row 113: Can't Eat Just One / Night Vinegar
$ws.Range("H113").Value = 1913.8
$ws.Range("J113").Value = 1913.8
$ws.Range("L113").Value = 5741.4
$ws.Range("N113").Value = -10081.4
# Row 117: A Good Omen / Peppered Popotoes
$ws.Range("H117").Value = 1504.75
$ws.Range("I117").Value = 1539.6666
$ws.Range("J117").Value = 1400
$ws.Range("K117").Value = 4618.9998
$ws.Range("L117").Value = 4200
$ws.Range("M117").Value = -1176.9998
$ws.Range("N117").Value = -11084
# Row 120: A Happy End / Paella
$ws.Range("H120").Value = 19490.666
$ws.Range("I120").Value = 9990
$ws.Range("K120").Value = 29970
$ws.Range("M120").Value = -25132
# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 4352.1304
$ws.Range("I131").Value = 4776.6665
$ws.Range("J131").Value = 4288.45
$ws.Range("K131").Value = 14329.9995
$ws.Range("L131").Value = 12865.35
$ws.Range("M131").Value = -9289.999500000002
$ws.Range("N131").Value = -22945.35
# Row 139: Najoothie / Wild Banana Blend
$ws.Range("H139").Value = 3125
$ws.Range("I139").Value = 3125
$ws.Range("K139").Value = 9375
$ws.Range("M139").Value = -4235

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers / Copper Ingot
$ws.Range("H2").Value = 367.6
$ws.Range("I2").Value = 54.142857
$ws.Range("K2").Value = 54.142857
$ws.Range("M2").Value = 58.857143
# Row 21: Forever 21K / Brass Ring
$ws.Range("H21").Value = 31000
$ws.Range("J21").Value = 31000
$ws.Range("L21").Value = 31000
$ws.Range("N21").Value = -31346
# Row 30: Dog Tags Are for Dogs / Brass Ring
$ws.Range("H30").Value = 31000
$ws.Range("J30").Value = 31000
$ws.Range("L30").Value = 31000
$ws.Range("N30").Value = -31210
# Row 52: It's My Business to Know Things / Red Coral Armillae
$ws.Range("H52").Value = 3000
$ws.Range("I52").Value = 3000
$ws.Range("K52").Value = 3000
$ws.Range("M52").Value = -2741
# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 6653.8857
$ws.Range("I126").Value = 6669.25
$ws.Range("K126").Value = 20007.75
$ws.Range("M126").Value = -17537.75

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore / Hard Leather
$ws.Range("H16").Value = 5891.4287
$ws.Range("I16").Value = 1180
$ws.Range("J16").Value = 15314.286
$ws.Range("K16").Value = 1180
$ws.Range("L16").Value = 15314.286
$ws.Range("M16").Value = -1010
$ws.Range("N16").Value = -15654.286
# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 4420
$ws.Range("J22").Value = 4833.3335
$ws.Range("L22").Value = 4833.3335
$ws.Range("N22").Value = -5423.3335
# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 4420
$ws.Range("J27").Value = 4833.3335
$ws.Range("L27").Value = 4833.3335
$ws.Range("N27").Value = -5047.3335
# Row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 1729.2142
$ws.Range("I93").Value = 1715.5834
$ws.Range("J93").Value = 1811
$ws.Range("K93").Value = 1715.5834
$ws.Range("L93").Value = 1811
$ws.Range("M93").Value = -467.5834
$ws.Range("N93").Value = -4307
# Row 104: Brace Yourselves / Gazelleskin Bracers of Fending
$ws.Range("H104").Value = 96919
$ws.Range("J104").Value = 96919
$ws.Range("L104").Value = 96919
$ws.Range("N104").Value = -103907
# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 7832.0605
$ws.Range("I132").Value = 7155.143
$ws.Range("K132").Value = 21465.429
$ws.Range("M132").Value = -18935.429

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107: Flax Wax / Bright Linen Yarn
$ws.Range("H107").Value = 1770.931
$ws.Range("I107").Value = 1430.8695
$ws.Range("K107").Value = 4292.6085
$ws.Range("M107").Value = -2372.6085
# Row 129: Lifetime of Gleaning / Scarlet Moko Beret of Gathering
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("M129").ClearContents()  # was -79990
$ws.Range("N129").ClearContents()  # was -100000
